# Replace "Melvin Hazen Valley Branch" lookup entries with a new
# "Reservation 630" entry across the lookup-table workbook.
#
# location_id: the id row keeps its id (TMH01) but now points at the new
#   display name.
# location_name: this sheet is keyed by name, so a brand-new row is
#   inserted (Reservation 630 / Reservation 630 / Piedmont) rather than
#   overwriting the Melvin Hazen Valley Branch row in place.
# connectivity_summary / fish_summary / habitat_summary /
#   macroinvertebrate_summary: each references the display name in
#   column A and needs the same rename.

$wb = $excel.ActiveWorkbook

# --- location_id: TMH01's name -> Reservation 630 -------------------------
$wsLocationId = $wb.Worksheets.Item("location_id")
$wsLocationId.Range("B19").Value = "Reservation 630"

# --- location_name: insert a new row for Reservation 630 ------------------
$wsLocationName = $wb.Worksheets.Item("location_name")
$wsLocationName.Rows.Item(19).Insert() | Out-Null
$wsLocationName.Range("A19").Value = "Reservation 630"
$wsLocationName.Range("B19").Value = "Reservation 630"
$wsLocationName.Range("C19").Value = "Piedmont"
# the old Melvin Hazen Valley Branch row (now row 18) now maps to the new name
$wsLocationName.Range("B18").Value = "Reservation 630"

# --- summary sheets: rename the row-16 entry -------------------------------
$wsConnectivity = $wb.Worksheets.Item("connectivity_summary")
$wsConnectivity.Range("A16").Value = "Reservation 630"

$wsFish = $wb.Worksheets.Item("fish_summary")
$wsFish.Range("A16").Value = "Reservation 630"

$wsHabitat = $wb.Worksheets.Item("habitat_summary")
$wsHabitat.Range("A16").Value = "Reservation 630"

$wsMacro = $wb.Worksheets.Item("macroinvertebrate_summary")
$wsMacro.Range("A16").Value = "Reservation 630"

# --- restore/update each sheet's last-used selection -----------------------
$wsLocationId.Range("A19").Select() | Out-Null

$wsConnectivity.Range("A16").Select() | Out-Null
$wsFish.Range("A16").Select() | Out-Null
$wsHabitat.Range("A16").Select() | Out-Null
$wsMacro.Range("A16").Select() | Out-Null

$wsEiaSubsheds = $wb.Worksheets.Item("eia_subsheds")
$wsEiaSubsheds.Range("C34").Select() | Out-Null

$wsEiaRockCreek = $wb.Worksheets.Item("eia_rock_creek")
$wsEiaRockCreek.Range("C8").Select() | Out-Null

$wsEia2023Temp = $wb.Worksheets.Item("eia_2023_TEMP")
$wsEia2023Temp.Range("P63").Select() | Out-Null

# location_name stays the active sheet/tab, so activate it last and leave
# its selection where the diff expects it.
$wsLocationName.Activate() | Out-Null
$wsLocationName.Range("O10").Select() | Out-Null
